$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header texts
$ws.Range("A1").Value = "Fund"
$ws.Range("B1").Value = "AUM"
$ws.Range("C1").Value = "M_Fee"
$ws.Range("D1").Value = "P_Fee"
$ws.Range("E1").Value = "Starting_Year"

# Fix fund fees (AUM values)
$ws.Range("B4").Value = 222
$ws.Range("B5").Value = 1500

# Header row is now shorter text, so it needs less height
$ws.Rows("1:1").RowHeight = 51

# Update selection to match target state
$ws.Range("B6").Select()
